# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.593.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.288.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '96.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.610'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0929'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.631.86'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.38'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.842'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.286.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.605.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  +1.60%  '
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("E22").Value = '  +12.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.10%  '
$ws.Range("E24").Value = '  -6.44%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.55'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("E28").Value = '  +2.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0897'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.36'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.08%  '
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("E36").Value = '  -2.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0352'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.35'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.241'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.55%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.46%  '
$ws.Range("E44").Value = '  +2.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.61%  '
$ws.Range("E47").Value = '  -0.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.185'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.73%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.513.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.66%  '
